# Word represents a manual line break (<w:br/>) as Chr(11) ("vertical tab")
# inside run text. Building replacement strings with [char]11 and assigning
# them straight to a Range's .Text (after using Find purely to *locate* the
# range, with no ReplaceWith) gives us literal <w:br/> splits without
# Word's Find/Replace "smart quotes" autocorrect mangling the apostrophes.
$lb = [char]11

$d = $word.ActiveDocument

# 1) "Identifiant(s) de(s) ressource(s)" -> "Identifiant(s) de la ressource(s) partagé(s)"
$rng = $d.Content
$rng.Find.Execute("Identifiant(s) de(s) ressource(s)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "Identifiant(s) de la ressource(s) partagé(s)"

# 2) "string" -> "string" + line break + "(REGEX: ...)"
$rng = $d.Content
$rng.Find.Execute("string", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "string" + $lb + "(REGEX: ^([\w-]+\.){3,4}resource(\.[\w-]+){1,2}$)"

# 3) "Liste des ID ... plus de détails" -> expanded multi-line description
$rng = $d.Content
$rng.Find.Execute("Liste des ID des ressources pour lesquels le demandeur a besoin d'obtenir plus de détails", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "Liste des ID des ressources pour lesquels le demandeur a besoin d'obtenir plus de détails. " + $lb + `
            "A valoriser avec l'identifiant partagé unique de la ressource engagée, normé comme suit :" + $lb + `
            "{orgID}.resource.{ID unique de la ressource partagée}" + $lb + `
            "OU - uniquement dans le cas où un ID unique de ressource ne peut pas être garanti par l'organisation propriétaire :" + $lb + `
            "{orgID}.resource.{sendercaseId}.{n° d’ordre chronologique de la ressource}"

# 4) "76_45101#SMUR1" -> two example identifiers separated by a line break
$rng = $d.Content
$rng.Find.Execute("76_45101#SMUR1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "fr.health.samu770.resource.VLM250" + $lb + "fr.health.samu440.resource.DRFR15DDXAAJJJ0000.1"
